# Refresh the hourly crypto-ranking snapshot (GitHub Actions bot update).
#
# Most "Price" (column D) values look like plain decimals (e.g. "10.00"),
# but the sheet stores them as TEXT, not numbers, so that formatting such
# as trailing zeros and thousands-style dot grouping (e.g. "33.781.62")
# survives untouched. Excel's Range.Value setter auto-detects plain
# numeric-looking strings and silently coerces them to real numbers
# (dropping trailing zeros etc.), so for any new price that would
# otherwise be mis-detected as a number we prefix it with a leading
# apostrophe to force a literal-text entry, exactly as a person typing
# the value into Excel would do to keep it as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Ref,
        [string]$Text
    )
    # Excel's Range.Value setter auto-coerces plain decimal-looking strings
    # into real numbers (losing e.g. trailing zeros). Detect that case with
    # a simple regex (no .NET parsing helpers available in this sandbox)
    # and force literal-text entry via a leading apostrophe, same as a
    # person typing it into the cell would do to keep it text.
    if ($Text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($Ref).Value = "'" + $Text
    } else {
        $ws.Range($Ref).Value = $Text
    }
}

# row 2 - Bitcoin
Set-TextValue "D2" "33.781.62"
Set-TextValue "E2" "  +8.32%  "

# row 3 - Ethereum
Set-TextValue "D3" "1.776.87"
Set-TextValue "E3" "  +4.45%  "

# row 4 - TetherUSD
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.09%  "

# row 5 - BNB
Set-TextValue "D5" "225.39"
Set-TextValue "E5" "  +1.86%  "

# row 6 - XRP
Set-TextValue "D6" "0.561"
Set-TextValue "E6" "  +4.99%  "

# row 7 - USDC (price unchanged)
Set-TextValue "E7" "  -0.09%  "

# row 8 - Solana
Set-TextValue "D8" "30.53"
Set-TextValue "E8" "  +2.74%  "

# row 9 - OKB
Set-TextValue "D9" "46.53"
Set-TextValue "E9" "  +2.70%  "

# row 10 - Cardano (price unchanged)
Set-TextValue "E10" "  +3.69%  "

# row 11 - Dogecoin
Set-TextValue "D11" "0.0665"
Set-TextValue "E11" "  +3.47%  "

# row 12 - TRON (price unchanged)
Set-TextValue "E12" "  +1.25%  "

# row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.028.24"
Set-TextValue "E13" "  +4.07%  "

# row 14 - WrappedEther
Set-TextValue "D14" "1.774.70"
Set-TextValue "E14" "  +4.16%  "

# row 15 - Polygon (price unchanged)
Set-TextValue "E15" "  +2.17%  "

# row 16 - WrappedBTC
Set-TextValue "D16" "33.760.41"
Set-TextValue "E16" "  +8.27%  "

# row 17 - Chainlink
Set-TextValue "D17" "10.00"
Set-TextValue "E17" "  -2.00%  "

# row 18 - Polkadot (price unchanged)
Set-TextValue "E18" "  +0.62%  "

# row 19 - Litecoin
Set-TextValue "D19" "68.49"
Set-TextValue "E19" "  +2.06%  "

# row 20 - BitcoinCash (volume unchanged)
Set-TextValue "D20" "252.06"

# row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0740"
Set-TextValue "E21" "  +2.24%  "

# row 22 - Dai (price unchanged)
Set-TextValue "E22" "  -0.02%  "

# row 23 - Avalanche
Set-TextValue "D23" "10.28"
Set-TextValue "E23" "  +1.47%  "

# row 24 - Uniswap (price unchanged)
Set-TextValue "E24" "  -2.22%  "

# row 25 - Toncoin (price unchanged)
Set-TextValue "E25" "  -1.08%  "

# row 26 - Monero
Set-TextValue "D26" "159.07"
Set-TextValue "E26" "  +0.18%  "

# row 27 - EthereumClassic
Set-TextValue "D27" "16.50"
Set-TextValue "E27" "  +3.23%  "

# row 28 - Stellar (price unchanged)
Set-TextValue "E28" "  +1.39%  "

# row 29 - Cosmos
Set-TextValue "D29" "6.94"
Set-TextValue "E29" "  +2.54%  "

# row 30 - BinanceUSD
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  -0.11%  "

# row 31 - Filecoin
Set-TextValue "D31" "3.82"
Set-TextValue "E31" "  +3.71%  "

# row 32 - Hedera
Set-TextValue "D32" "0.0513"
Set-TextValue "E32" "  +1.95%  "

# row 33 - PancakeSwap (price unchanged)
Set-TextValue "E33" "  +3.63%  "

# row 34 - InternetComputer(DFINITY) (price unchanged)
Set-TextValue "E34" "  +5.18%  "

# row 35 - LidoDAOToken
Set-TextValue "D35" "1.84"
Set-TextValue "E35" "  +5.82%  "

# row 36 - Maker
Set-TextValue "D36" "1.483.50"
Set-TextValue "E36" "  -1.95%  "

# row 37 - TrustWalletToken (price unchanged)
Set-TextValue "E37" "  +3.13%  "

# row 38 - ImmutableX (price unchanged)
Set-TextValue "E38" "  +3.12%  "

# rows 39/40 - Aave and VeChain swapped ranking order
Set-TextValue "B39" "Aave"
Set-TextValue "C39" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D39" "83.28"
Set-TextValue "E39" "  -0.05%  "

Set-TextValue "B40" "VeChain"
Set-TextValue "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.0185"
Set-TextValue "E40" "  +2.36%  "

# row 41 - HuobiToken (price unchanged)
Set-TextValue "E41" "  +1.66%  "

# row 42 - MXToken
Set-TextValue "D42" "2.70"
Set-TextValue "E42" "  +0.05%  "

# row 43 - ARBITRUM
Set-TextValue "D43" "0.885"
Set-TextValue "E43" "  +4.06%  "

# row 44 - RenderToken (price unchanged)
Set-TextValue "E44" "  +2.25%  "

# row 45 - Kaspa
Set-TextValue "D45" "0.0513"
Set-TextValue "E45" "  +1.65%  "

# row 46 - WEMIXToken (price unchanged)
Set-TextValue "E46" "  +3.34%  "

# row 47 - RocketPoolETH
Set-TextValue "D47" "1.927.82"
Set-TextValue "E47" "  +4.71%  "

# row 48 - FraxShare
Set-TextValue "D48" "5.72"
Set-TextValue "E48" "  +2.57%  "

# row 49 - PaxDollar
Set-TextValue "D49" "0.999"
Set-TextValue "E49" "  -0.14%  "

# row 50 - InjectiveProtocol
Set-TextValue "D50" "11.78"
Set-TextValue "E50" "  +13.92%  "

# row 51 - BitcoinSV
Set-TextValue "D51" "50.73"
Set-TextValue "E51" "  -2.88%  "
